$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Resize / reposition the two "Conector Reto" straight-line connectors
#    (Shape geometry is expressed in points through the COM object model;
#    the target EMU values from the authoring tool convert at 12700 EMU/pt)
# ---------------------------------------------------------------------------
$s5 = $d.Shapes.Item("Conector Reto 5")
$s5.Left   = 45720 / 12700
$s5.Width  = 25400 / 12700
$s5.Height = 3500120 / 12700

$s9 = $d.Shapes.Item("Conector Reto 9")
$s9.Width  = 5080 / 12700
$s9.Height = 800100 / 12700

# ---------------------------------------------------------------------------
# 2) Text edits
# ---------------------------------------------------------------------------

$d.Content.Find.Execute("Framework. APIs", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Framework e Angular 1.x. APIs", 2)

$d.Content.Find.Execute("AWS Lambdas, Containers, testes com XUnit e Moq, RabbitMq, SQS", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "AWS Lambdas, containers, testes com XUnit, ferramentas de mensagerias como rabbitMq e SQS", 2)

$d.Content.Find.Execute("com grande impacto de negócio.", $true, $false, $false, $false, $false, `
    $true, 1, $false, "com grande impacto de negócio. Utilzei o ASP.NET Razor Pages para a construção dessa ferramenta.", 2)

$d.Content.Find.Execute("criação de scripts em python, ", $true, $false, $false, $false, $false, `
    $true, 1, $false, "criação de scripts em python e ", 2)

# ---------------------------------------------------------------------------
# 3) Split "TCC realizado" into "TCC realiz" + "ado ..." and relocate the
#    "_GoBack" bookmark to the split point (mirrors where Word last left the
#    editing cursor after this change).
# ---------------------------------------------------------------------------
$tcc = $d.Content
$tcc.Find.Execute("TCC realizado com o tema sobre desenvolvimento web: ")
$splitPoint = $tcc.Start + 10

$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

$bmRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------------
# 4) Mark the built-in "footer" style as a Quick Style (w:qFormat)
# ---------------------------------------------------------------------------
$footerStyle = $d.Styles.Item("footer")
$footerStyle.QuickStyle = $true

Write-Host "done"
